# The workbook opens with "Simplified contract multiplier" as the active
# sheet (it matches $wb.ActiveSheet), which is where the new asset row goes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 ("kappa") is the last populated row (A1:G15). Carry its formatting
# down onto the new row 16 before writing the new asset's data into it, so
# the label cell (A16) picks up the same font/alignment style as A15.
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null

# New commodity row: label in column A, and the two factor columns F/G
# (same columns used by every other row) filled in with the new values.
$ws.Range("A16").Value = "fake_asset"
$ws.Range("F16").Value = 0.1
$ws.Range("G16").Value = 0.01

# Leave the selection where the author's session ended up.
$ws.Range("C11").Select() | Out-Null
